# "add support for mac m1"
# Update TestSteps sheet: change row 4 (TS003) and row 5 (TS004) content,
# and remove the trailing rows 6-8 (TS005/TS006/TS007) which are no longer used.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestSteps")

# Row 4 (TS003): input username -> input info / hello / input
$ws.Range("C4").Value = "输入信息"
$ws.Range("D4").Value = "hello"
$ws.Range("E4").Value = "input"
$ws.Range("F4").Value = "input"

# Row 5 (TS004): input password -> search
$ws.Range("C5").Value = "搜索"
$ws.Range("D5").ClearContents()
$ws.Range("E5").Value = "search"
$ws.Range("F5").Value = "click"

# Rows 6-8 (TS005, TS006, TS007) are dropped entirely
$ws.Range("A6:G8").ClearContents()

# New selected cell in the saved view
$ws.Range("E5").Select()
